{"js": "// Highlight the four checklist paragraphs discussing the Twillio/phone-controller\n// work in yellow, and (while we're re-doing this section) collapse the three\n// runs that spell out \"Add campa\" + \"i\" + \"gn category...\" back into a single\n// run of contiguous text \u2014 matching the canonical edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Identify the target paragraphs by their distinctive leading text rather than\n// a hard-coded index, so the script is resilient to minor document changes:\n//   \"Finish connecting callee details to phone controller view...\"\n//   \"Add campaign category in model and have campaign displayed in details\"\n//   \"Have \\u201cdial\\u201d and \\u201ctext\\u201d views triggered on button...\"\n//   \"Research twillio/stackof to display callee text responses...\"\nconst highlightStartsWith = [\n  \"Finish connecting\",\n  \"Add campa\",\n  \"Have \\u201cdial\\u201d and \\u201ctext\\u201d views triggered\",\n  \"Research \",\n];\n\nlet addCampaignParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  const isTarget = highlightStartsWith.some((snippet) => text.indexOf(snippet) === 0);\n  if (isTarget) {\n    paragraph.font.highlightColor = \"yellow\";\n  }\n  if (text.indexOf(\"Add campa\") === 0) {\n    addCampaignParagraph = paragraph;\n  }\n}\nawait context.sync();\n\n// The \"Add campaign...\" paragraph originally had its text split across three\n// runs (\"Add campa\" / \"i\" / \"gn category in model and have campaign displayed\n// in \"). Re-insert the combined text as a single run (this also keeps the\n// yellow highlight already applied above).\nif (addCampaignParagraph !== null) {\n  const mergedTextRanges = addCampaignParagraph.search(\n    \"Add campaign category in model and have campaign displayed in \",\n    { matchCase: true }\n  );\n  await context.sync();\n\n  if (mergedTextRanges.items.length > 0) {\n    mergedTextRanges.items[0].insertText(\n      \"Add campaign category in model and have campaign displayed in \",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight the four checklist paragraphs discussing the Twillio/phone-controller\n# work in yellow, and (while we're re-doing this section) collapse the three\n# runs that spell out \"Add campa\" + \"i\" + \"gn category...\" back into a single\n# run of contiguous text -- matching the canonical edit.\n\n$d = $word.ActiveDocument\n$wdYellow = 7\n\n# Identify the target paragraphs by their distinctive leading text rather than\n# a hard-coded index, so the script is resilient to minor document changes:\n#   \"Finish connecting callee details to phone controller view...\"\n#   \"Add campaign category in model and have campaign displayed in details\"\n#   \"Have \"dial\" and \"text\" views triggered on button...\"\n#   \"Research twillio/stackof to display callee text responses...\"\n$highlightPrefixes = @(\n    \"Finish connecting\",\n    \"Add campa\",\n    \"Have \",\n    \"Research \"\n)\n\n$paragraphCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $paragraphCount; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n    foreach ($prefix in $highlightPrefixes) {\n        if ($text.StartsWith($prefix)) {\n            # Set via .Font (not the bare Range) so the paragraph-mark run\n            # properties (w:pPr/w:rPr) pick up the highlight too, not just\n            # the visible runs.\n            $p.Range.Font.HighlightColorIndex = $wdYellow\n            break\n        }\n    }\n}\n\n# The \"Add campaign...\" paragraph originally had its text split across three\n# runs (\"Add campa\" / \"i\" / \"gn category in model and have campaign displayed\n# in \"). Re-insert the combined text as a single run (this also keeps the\n# yellow highlight already applied above) via Find & Replace.\n$mergeRange = $d.Range()\n$mergeRange.Find.Execute(\n    \"Add campaign category in model and have campaign displayed in \",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"Add campaign category in model and have campaign displayed in \",\n    2\n)\n"}
